$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K: header + per-row values (mirrors the pattern of a new
# package release column being appended to the version matrix). Row 5 is
# intentionally left untouched - it did not receive a K cell in the source
# edit either.
$ws.Range("K1").Value = "12.1.0"
$ws.Range("K2").Value = "10.1.0"
$ws.Range("K4").Value = "3.4.2"
$ws.Range("K8").Value = "9.0.1"

# Copy formatting from column J (the previous last column) onto column K so
# the new column matches the existing per-row styling (bold header, italic
# J3-style cell, plain body rows, etc.) Skip row 5 so it keeps no K cell.
$ws.Range("J1:J4").Copy() | Out-Null
$ws.Range("K1:K4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("J6:J11").Copy() | Out-Null
$ws.Range("K6:K11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Move the active selection to K1, matching the post-edit workbook state.
$ws.Range("K1").Select() | Out-Null
